$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the two blank rows (2:3) so the test rows move up ---
$ws.Rows("2:3").Delete()

# --- Step 2: shift the old C:F block two columns to the right (E:H) ---
$ws.Columns("A:B").Insert()

# --- Step 3: drop the now-empty A:D block, bringing old C->A, D->B, E->C, F->D ---
$ws.Columns("A:D").Delete()

# Now: A = "Prueba" column (width 103), B = empty spacer (width 14),
#      C = "Resultado" column (width 18.85546875), D = empty spacer (width 11.42578125)
#      holding the lone formatted marker cell at D1, with A10/B10 both carrying that marker style.

# --- Step 4: clear the leftover marker cells/content that must not survive the move ---
$ws.Range("D1").Clear()
$ws.Range("B10").Clear()

# --- Step 5: drop the now fully empty spacer column D entirely ---
$ws.Columns("D").Delete()

# --- Step 6: build the new "Resultado Esperado" column in B ---
$ws.Range("B1").Value = "Resultado Esperado"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = "Si"
    $ws.Cells.Item($r, 3).Value = "Si"
    $ws.Cells.Item($r, 3).Style = "Good"
}
$ws.Columns("B").ColumnWidth = 20.7

# --- Step 7: recreate the two lone formatted marker cells at their new spots ---
$ws.Range("E1").Font.Underline = $true
$ws.Range("C19").Font.Underline = $true

# --- Step 8: selection matches the authored workbook ---
$ws.Range("C19").Select()
